$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 62, shifting existing rows 62-75 down to 63-76
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new data record
$ws.Cells.Item(62, 1).Value = 9
$ws.Cells.Item(62, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44641
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = 100114002
$ws.Cells.Item(62, 7).Value = "Camote"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 1600
$ws.Cells.Item(62, 11).Value = 10000
$ws.Cells.Item(62, 12).Value = 11000
$ws.Cells.Item(62, 13).Value = 10500
$ws.Cells.Item(62, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(62, 15).Value = "Perú"
$ws.Cells.Item(62, 16).Value = 583
$ws.Cells.Item(62, 17).Value = 18
$ws.Cells.Item(62, 18).Value = "Hortaliza"
